$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$images = @("flower/flower096.png","flower/flower073.png","face/face067.png","flower/flower112.png","flower/flower093.png","face/face105.png","flower/flower083.png","flower/flower072.png","face/face073.png","face/face110.png","flower/flower105.png","face/face074.png","flower/flower068.png","face/face080.png","face/face089.png","flower/flower076.png","flower/flower067.png","face/face081.png","face/face090.png","flower/flower109.png","flower/flower104.png","face/face075.png","flower/flower098.png","flower/flower123.png","face/face099.png","flower/flower120.png","flower/flower100.png","face/face069.png","face/face102.png","face/face071.png","face/face072.png","face/face103.png")
$words = @("töten","haken","antun","währen","sondern","stärken","hauen","spielen","loben","schenken","runden","schicken","fühlen","raten","wiegen","füttern","drehen","strahlen","segeln","drohen","hoffen","starten","gründen","rasen","ehren","rücken","mieten","schätzen","tagen","husten","lehnen","kehren")
$cats = @("flower","flower","face","flower","flower","face","flower","flower","face","face","flower","face","flower","face","face","flower","flower","face","face","flower","flower","face","flower","flower","face","flower","flower","face","face","face","face","face")

for ($i = 0; $i -lt 32; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $images[$i]
    $ws.Cells.Item($row, 3).Value = $words[$i]
    $ws.Cells.Item($row, 4).Value = $cats[$i]
}
